# Small typo fix: the department now has two AWS DeepRacer vehicles,
# not just one - update the "Ausgangslage" paragraph accordingly.
#
#   "Der Abteilung Informatik steht nun ein AWS DeepRacer zur Verfügung."
# becomes
#   "Der Abteilung Informatik stehen nun zwei AWS DeepRacer zur Verfügung."

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Der Abteilung Informatik steht nun ein AWS DeepRacer zur Verfügung.",
    $false,
    $true,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Der Abteilung Informatik stehen nun zwei AWS DeepRacer zur Verfügung.",
    2
)
